# Applies the edit described by the commit:
# "agregado suma de todos los sectores al terminar los resultados de
#  potencia por sectores" — marks the sector-power spec rows as
# completed ("Listo") in the "Cálculos Sectores" sheet, and leaves that
# sheet as the active/selected one (mirroring the tab/selection switch
# seen in the diff, which moves from "Calculos Fs" to "Cálculos Sectores").

$wb = $excel.ActiveWorkbook

$wsFs = $wb.Worksheets.Item("Calculos Fs")
$wsSectores = $wb.Worksheets.Item("Cálculos Sectores")

# Mark the first five spec rows (ID 3010-3050) of the sectors sheet as
# finished, in the "Estado" (status) column.
$wsSectores.Range("C2:C6").Value = "Listo"

# Leave "Calculos Fs" scrolled down a bit (row 4 at the top) with its
# previous selection (B6), just no longer the active tab.
$wsFs.Activate()
$wsFs.Range("B6").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

# Switch to the "Cálculos Sectores" sheet, scroll it up to row 4 and
# select C6, matching the new selection left behind after filling in
# the status column.
$wsSectores.Activate()
$wsSectores.Range("C6").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
